# today's client changes 2
#
# Replaces the sample candidate rows (rows 2-5) with a new batch of
# candidates, and adjusts the sheet view / a couple of column widths
# that Excel re-derived after the data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 2 : Dinesh Kartick ----
$ws.Range("A2").Value = "Calling"
$ws.Range("B2").Value = "Facebook"
$ws.Range("C2").Value = "Subham Ghosh"
$ws.Range("D2").Value = "2022-10-10"
$ws.Range("E2").Value = "Dinesh Kartick"
$ws.Range("F2").Value = "Male"
$ws.Range("G2").Value = "1999-04-05"
$ws.Range("H2").Value = "Higer Secondary"
$ws.Range("I2").Value = "BCA"
$ws.Range("J2").Value = 7410258900
$ws.Range("K2").Value = 5242525555
$ws.Range("L2").Value = "dinesh@yopmail.com"
$ws.Range("M2").Value = 917485758585
$ws.Range("N2").Value = "Kokata "
$ws.Range("O2").Value = "Hindu"
$ws.Range("P2").Value = "ECR"
$ws.Range("Q2").Value = "Four Wheeler"
$ws.Range("R2").Value = "Two Wheeler"
$ws.Range("S2").Value = "Poor"
$ws.Range("T2").Value = "Poor"
$ws.Range("U2").Value = "S5632328"
$ws.Range("V2").Value = "Yes"
$ws.Range("W2").Value = "MACHINE EMBROIDER"
$ws.Range("X2").Value = "OFFICE MACHINE OPERATOR"
$ws.Range("Y2").Value = "FINISHING CARPENTER"
$ws.Range("Z2").Value = "1 Year Experience"
$ws.Range("AA2").Value = "1 Year Experience"
$ws.Range("AB2").Value = "Good Work"

# ---- Row 3 : Bipul Das ----
$ws.Range("A3").Value = "Calling"
$ws.Range("B3").Value = "Instagram"
$ws.Range("C3").Value = "Azhar SK"
$ws.Range("D3").Value = "2022-04-04"
$ws.Range("E3").Value = "Bipul Das"
$ws.Range("F3").Value = "Male"
$ws.Range("G3").Value = "1998-10-10"
$ws.Range("H3").Value = "MBA"
$ws.Range("I3").Value = "BSC"
$ws.Range("J3").Value = 8527418520
$ws.Range("K3").Value = 5242525555
$ws.Range("L3").Value = "bipul@yopmail.com"
$ws.Range("M3").Value = 917485857585
$ws.Range("N3").Value = "Kokata "
$ws.Range("O3").Value = "Hindu"
$ws.Range("P3").Value = "ENCR"
$ws.Range("Q3").Value = "Four Wheeler"
$ws.Range("R3").Value = "Two Wheeler"
$ws.Range("S3").Value = "Basic"
$ws.Range("T3").Value = "Basic"
$ws.Range("U3").Value = "P7418525"
$ws.Range("V3").Value = "Yes"
$ws.Range("W3").Value = "MECHANICAL SUPERVISOR"
$ws.Range("X3").Value = "LMV MECHANIC PETROL"
$ws.Range("Y3").Value = "GARDENER"
$ws.Range("Z3").Value = "1 Year Experience"
$ws.Range("AA3").Value = "1 Year Experience"
$ws.Range("AB3").Value = "Good Work"

# ---- Row 4 : Nirmal Ghosh ----
$ws.Range("A4").Value = "Calling"
$ws.Range("B4").Value = "Telecalling"
$ws.Range("C4").Value = "Dilip Ghosh"
$ws.Range("D4").Value = "2021-04-04"
$ws.Range("E4").Value = "Nirmal Ghosh"
$ws.Range("F4").Value = "Female"
$ws.Range("G4").Value = "1997-10-10"
$ws.Range("H4").Value = "Higer Secondary"
$ws.Range("I4").Value = "B.Tech"
$ws.Range("J4").Value = 9517539620
$ws.Range("K4").Value = 5242525555
$ws.Range("L4").Value = "nirmal@yopmail.com"
$ws.Range("M4").Value = 917485235695
$ws.Range("N4").Value = "Kokata "
$ws.Range("O4").Value = "Hindu"
$ws.Range("P4").Value = "ECR"
$ws.Range("Q4").Value = "Four Wheeler"
$ws.Range("R4").Value = "Two Wheeler"
$ws.Range("S4").Value = "No"
$ws.Range("T4").Value = "Good"
$ws.Range("U4").Value = "Y8575858"
$ws.Range("V4").Value = "No"
$ws.Range("W4").Value = "MOBILE CRANE OPERATOR"
$ws.Range("X4").Value = "AREA RESTURANT MANAGER"
$ws.Range("Y4").Value = "WELDER"
$ws.Range("Z4").Value = "1 Year Experience"
$ws.Range("AA4").Value = "1 Year Experience"
$ws.Range("AB4").Value = "Good Work"

# ---- Row 5 : Jahiralom Sk ----
$ws.Range("A5").Value = "Calling"
$ws.Range("B5").Value = "Facebook"
$ws.Range("C5").Value = "Azhar SK"
$ws.Range("D5").Value = "2022-10-10"
$ws.Range("E5").Value = "Jahiralom Sk"
$ws.Range("F5").Value = "Male"
$ws.Range("G5").Value = "1995-10-10"
$ws.Range("H5").Value = "BBA"
$ws.Range("I5").Value = "BSC"
$ws.Range("J5").Value = 8545632102
$ws.Range("K5").Value = 5242525555
$ws.Range("L5").Value = "jahiralom@yopmail.com"
$ws.Range("M5").Value = 917485230230
$ws.Range("N5").Value = "Kokata "
$ws.Range("O5").Value = "Muslim"
$ws.Range("P5").Value = "ENCR"
$ws.Range("Q5").Value = "Four Wheeler"
$ws.Range("R5").Value = "Two Wheeler"
$ws.Range("S5").Value = "Good"
$ws.Range("T5").Value = "No"
$ws.Range("U5").Value = "Z8575856"
$ws.Range("V5").Value = "No"
$ws.Range("W5").Value = "OFFICE BOY"
$ws.Range("X5").Value = "GLASS DESIGNER"
$ws.Range("Y5").Value = "GENTS TAILOR"
$ws.Range("Z5").Value = "1 Year Experience"
$ws.Range("AA5").Value = "1 Year Experience"
$ws.Range("AB5").Value = "Good Work"

# ---- Selection / scroll position ----
# was: topLeftCell N1, selection W9 -> now: default top-left, selection E10
[void]$ws.Range("E10").Select()

# ---- Column width tweaks (Excel re-derived these after the content changed) ----
$ws.Columns.Item(13).ColumnWidth = 21
$ws.Columns.Item(23).ColumnWidth = 23.6
